# Refresh cryptocurrency price/volume data, preserving original text-cell formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.863.66'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -7.61%  '
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.651.37'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -9.24%  '
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.020'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +1.49%  '
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '221.14'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -5.11%  '
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.032'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.57%  '
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5046'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -14.14%  '
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2574'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -5.50%  '
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.78'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.02%  '
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06125'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -9.64%  '
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07466'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.62%  '
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.723.78'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -5.04%  '
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.468'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.66%  '
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'Polygon'
$ws.Range("B14").ClearFormats()
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("C14").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5723'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -7.85%  '
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("B15").ClearFormats()
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("C15").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.846.78'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -10.57%  '
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008042'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -12.82%  '
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("B17").ClearFormats()
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("C17").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.51'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -13.16%  '
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("B18").ClearFormats()
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("C18").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.648.52'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -7.58%  '
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("B19").ClearFormats()
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("C19").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.033'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -6.78%  '
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'Dai'
$ws.Range("B20").ClearFormats()
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("C20").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.015'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.84'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.36%  '
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '185.26'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -10.51%  '
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.015'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.79%  '
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.160'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -8.72%  '
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.04'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -6.41%  '
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.638'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.42%  '
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1155'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -8.03%  '
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.15'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -6.37%  '
$ws.Range("E28").ClearFormats()

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.360'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.45%  '
$ws.Range("E29").ClearFormats()

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05821'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -7.41%  '
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.365'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.50%  '
$ws.Range("E31").ClearFormats()

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.473'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.85%  '
$ws.Range("E32").ClearFormats()

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.384'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -8.11%  '
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.571'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -6.30%  '
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("B35").ClearFormats()
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("C35").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9831'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.43%  '
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("B36").ClearFormats()
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("C36").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.477'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.31%  '
$ws.Range("E36").ClearFormats()

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("B37").ClearFormats()
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C37").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6006'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.77%  '
$ws.Range("E37").ClearFormats()

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'MXToken'
$ws.Range("B38").ClearFormats()
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("C38").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.593'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -5.65%  '
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01587'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -6.33%  '
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8707'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.95%  '
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.035'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.84%  '
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.064.96'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.47%  '
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.856'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -8.33%  '
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("B44").ClearFormats()
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("C44").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.857.77'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.81%  '
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Quant'
$ws.Range("B45").ClearFormats()
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("C45").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '94.53'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -5.18%  '
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4424'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.30%  '
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Aave'
$ws.Range("B48").ClearFormats()
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("C48").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.04'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -8.02%  '
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Cronos'
$ws.Range("B49").ClearFormats()
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("C49").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05239'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.45%  '
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Frax'
$ws.Range("B50").ClearFormats()
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("C50").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9507'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -6.25%  '
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'USDD'
$ws.Range("B51").ClearFormats()
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("C51").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.023'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.28%  '
$ws.Range("E51").ClearFormats()

